$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.514.80"
$ws.Range("E2").Value = "  +0.12%  "

$ws.Range("D3").Value = "2.301.17"
$ws.Range("E3").Value = "  -0.17%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'316.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.57%  "

$ws.Range("D6").Value = "'103.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.18%  "

$ws.Range("E7").Value = "  -0.38%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  -0.85%  "

$ws.Range("D10").Value = "'39.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.36%  "

$ws.Range("D11").Value = "'0.0904"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.70%  "

$ws.Range("E12").Value = "  +2.39%  "

$ws.Range("E13").Value = "  +1.43%  "

$ws.Range("E14").Value = "  +3.71%  "

$ws.Range("D15").Value = "'15.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.17%  "

$ws.Range("D16").Value = "2.649.56"
$ws.Range("E16").Value = "  -0.02%  "

$ws.Range("D17").Value = "2.307.00"
$ws.Range("E17").Value = "  +0.30%  "

$ws.Range("D18").Value = "42.606.10"
$ws.Range("E18").Value = "  +0.54%  "

$ws.Range("E19").Value = "  +1.39%  "

$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "'13.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +24.06%  "

$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "'0.0000106"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.44%  "

$ws.Range("E22").Value = "  +0.57%  "

$ws.Range("D23").Value = "'3.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.24%  "

$ws.Range("D24").Value = "'266.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.76%  "

$ws.Range("E25").Value = "  -1.48%  "

$ws.Range("D26").Value = "'1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.40%  "

$ws.Range("D27").Value = "'10.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.37%  "

$ws.Range("E28").Value = "  -1.69%  "

$ws.Range("D29").Value = "'22.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.69%  "

$ws.Range("E30").Value = "  +11.98%  "

$ws.Range("D31").Value = "'37.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.62%  "

$ws.Range("D32").Value = "'165.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.46%  "

$ws.Range("E33").Value = "  +0.88%  "

$ws.Range("E34").Value = "  -2.97%  "

$ws.Range("E35").Value = "  -4.73%  "

$ws.Range("E36").Value = "  -2.02%  "

$ws.Range("E37").Value = "  -1.31%  "

$ws.Range("E38").Value = "  +1.39%  "

$ws.Range("D39").Value = "'3.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.31%  "

$ws.Range("E40").Value = "  -1.70%  "

$ws.Range("E41").Value = "  +9.23%  "

$ws.Range("D42").Value = "'70.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.32%  "

$ws.Range("D43").Value = "'95.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.19%  "

$ws.Range("E44").Value = "  +0.82%  "

$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("D46").Value = "'12.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.48%  "

$ws.Range("D47").Value = "'117.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.47%  "

$ws.Range("D48").Value = "'80.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.75%  "

$ws.Range("D49").Value = "1.661.57"
$ws.Range("E49").Value = "  +3.51%  "

$ws.Range("D50").Value = "'5.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.34%  "

$ws.Range("D51").Value = "'8.85"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.73%  "
